$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.561.27"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.259.08"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.90"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.633"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.85"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.628"
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.61"
$ws.Range("E10").Value = "  +8.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0955"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.33"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.596.40"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.70"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.861"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.254.02"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.354.84"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("E19").Value = "  +3.51%  "
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.51"
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.24"
$ws.Range("E22").Value = "  +4.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.48"
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.16"
$ws.Range("E24").Value = "  +36.73%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.49"
$ws.Range("E26").Value = "  +2.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.62"
$ws.Range("E27").Value = "  -4.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.33"
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  +1.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.84"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.78"
$ws.Range("E31").Value = "  +0.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0835"
$ws.Range("E32").Value = "  -1.49%  "
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.87"
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.34"
$ws.Range("E35").Value = "  +13.66%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.54"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("E38").Value = "  +6.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.83"
$ws.Range("E39").Value = "  +4.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.20"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.85"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "63.42"
$ws.Range("E42").Value = "  +4.96%  "
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "108.95"
$ws.Range("E44").Value = "  -8.97%  "
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.998"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.20"
$ws.Range("E48").Value = "  +2.41%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("E50").Value = "  +5.38%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.16"
$ws.Range("E51").Value = "  -0.47%  "
